# Update "Pais" sheet: refresh the COVID country/provincia stats snapshot.
#
# The source data is a daily ranking table (countries sorted by total
# cases). Re-running the scrape changed several countries' totals, which in
# turn reshuffled a few rows that were close in rank to their neighbours.
# The net effect for this snapshot is: row-by-row numeric refreshes, plus
# three row-pairs where the country name at a fixed row flips with its
# neighbour (because the rank order of those two countries swapped while
# everything else's position stayed put), plus a refreshed "as of" timestamp
# string at the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Title row -------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 6 de Septiembre de 2020 a las 01:46"

# --- Data rows (row -> Pais, Casos totales, Nuevos casos, Casos activos,
#                Recuperados, Casos criticos, Muertes hoy, Muertes) -------
$rows = @{
    4   = @("Estados Unidos",        6428993, 39936, 3705937, 2530242, 0, 703, 192814)
    5   = @("Brasil",                4123000, 31199, 3296702,  700095, 0, 619, 126203)
    13  = @("Argentina",              471806,  9924,  340381,  121686, 0, 116,   9739)
    36  = @("Panama",                  96305,   709,   69223,   25007, 0,  12,   2075)
    54  = @("Nigeria",                 54905,   162,   42922,   10929, 0,   3,   1054)
    55  = @("Barein",                  54771,   676,   50645,    3930, 0,   1,    196)
    72  = @("Chequia",                 27752,   503,   19039,    8282, 0,   2,    431)
    84  = @("Bulgaria",                17050,    96,   12132,    4247, 0,   6,    671)
    92  = @("Noruega",                 11296,    65,    9348,    1684, 0,   0,    264)
    102 = @("Haiti",                    8336,    10,    5933,    2191, 0,   0,    212)
    106 = @("Luxemburgo",               6896,    42,    6126,     646, 0,   0,    124)
    109 = @("Montenegro",               5422,   147,    4224,    1091, 0,   1,    107)
    110 = @("Republica de Yibuti",      5387,     0,    5327,       0, 0,   0,     60)
    134 = @("Angola",                   2935,    59,    1192,    1626, 0,   2,    117)
    142 = @("Trinidad yTobago",         2230,   190,     717,    1480, 0,   2,     33)
    143 = @("Benin",                    2194,     0,    1793,     361, 0,   0,     40)
    151 = @("Uruguay",                  1669,    16,    1459,     165, 0,   0,     45)
    167 = @("Santo Tome y Principe",     898,     1,     859,      24, 0,   0,     15)
    169 = @("San Marino",                716,     1,     660,      14, 0,   0,     42)
    189 = @("Barbados",                  178,     0,     154,      17, 0,   0,      7)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = $vals[6]
    $ws.Cells.Item($r, 8).Value = $vals[7]
}
